$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 (H2:K2) - I2 and J2 change from text to numeric values
$ws.Range("H2").Value = 4.7830000000000004
$ws.Range("I2").Value = 9.5419999999999998
$ws.Range("J2").Value = 0.81
$ws.Range("K2").Value = 2.3740000000000001

# Update row 3 (H3:K3)
$ws.Range("H3").Value = 2.41
$ws.Range("I3").Value = 4.827
$ws.Range("J3").Value = 0.44400000000000001
$ws.Range("K3").Value = 1.2210000000000001

# Update row 4 (H4:K4)
$ws.Range("H4").Value = 1.2250000000000001
$ws.Range("I4").Value = 2.4249999999999998
$ws.Range("J4").Value = 0.22500000000000001
$ws.Range("K4").Value = 0.622

# Update row 5 (H5:K5)
$ws.Range("H5").Value = 0.627
$ws.Range("I5").Value = 1.1930000000000001
$ws.Range("J5").Value = 0.13200000000000001
$ws.Range("K5").Value = 0.32900000000000001

# Update selection
$ws.Range("O15").Select()
